$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KNX Group Addresses")
$ws.Name = "KNX GAs"
